$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 99.375541
$ws.Range("H2").Value = 298.126623
$ws.Range("I2").Value = 0.0220419662695056
$ws.Range("J2").Value = 0.0220419662695056
$ws.Range("M2").Value = 5.097982
$ws.Range("N2").Value = 15.293946
$ws.Range("O2").Value = 0.1960698660397332
$ws.Range("P2").Value = 0.1960698660397332
$ws.Range("Q2").Value = 506.614719258262
$ws.Range("R2").Value = 4559.532473324358
$ws.Range("S2").Value = 0.00432176537371428
$ws.Range("T2").Value = 0.00432176537371428

# Row 3
$ws.Range("G3").Value = 99.375541
$ws.Range("H3").Value = 298.126623
$ws.Range("I3").Value = 0.0220419662695056
$ws.Range("J3").Value = 0.0220419662695056
$ws.Range("O3").Value = 0.5754666399803534
$ws.Range("P3").Value = 0.5754666399803534
$ws.Range("Q3").Value = 1486.918291651518
$ws.Range("R3").Value = 13382.26462486366
$ws.Range("S3").Value = 0.01268441626767267
$ws.Range("T3").Value = 0.01268441626767267

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 99.375541
$ws.Range("H4").Value = 298.126623
$ws.Range("I4").Value = 0.0220419662695056
$ws.Range("J4").Value = 0.0220419662695056
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.940243666666667
$ws.Range("N4").Value = 17.820731
$ws.Range("O4").Value = 0.2284634939799134
$ws.Range("P4").Value = 0.2284634939799134
$ws.Range("Q4").Value = 590.3149280468238
$ws.Range("R4").Value = 5312.834352421413
$ws.Range("S4").Value = 0.005035784628118647
$ws.Range("T4").Value = 0.005035784628118647

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 4345.262939666666
$ws.Range("H5").Value = 13035.788819
$ws.Range("I5").Value = 0.9637999268679746
$ws.Range("J5").Value = 0.9637999268679746
$ws.Range("M5").Value = 5.097982
$ws.Range("N5").Value = 15.293946
$ws.Range("O5").Value = 0.1960698660397332
$ws.Range("P5").Value = 0.1960698660397332
$ws.Range("Q5").Value = 22152.07225168775
$ws.Range("R5").Value = 199368.6502651898
$ws.Range("S5").Value = 0.1889721225501084
$ws.Range("T5").Value = 0.1889721225501084

# Row 6
$ws.Range("D6").Value = "FAPs"
$ws.Range("I6").Value = 0.9637999268679746
$ws.Range("J6").Value = 0.9637999268679746
$ws.Range("M6").Value = 14.96261833333333
$ws.Range("N6").Value = 44.887855
$ws.Range("O6").Value = 0.5754666399803534
$ws.Range("P6").Value = 0.5754666399803534
$ws.Range("Q6").Value = 65016.51092421035
$ws.Range("R6").Value = 585148.5983178932
$ws.Range("S6").Value = 0.5546347055280236
$ws.Range("T6").Value = 0.5546347055280236

# Row 7
$ws.Range("D7").Value = "MuSCs"
$ws.Range("I7").Value = 0.9637999268679746
$ws.Range("J7").Value = 0.9637999268679746
$ws.Range("M7").Value = 5.940243666666667
$ws.Range("N7").Value = 17.820731
$ws.Range("O7").Value = 0.2284634939799134
$ws.Range("P7").Value = 0.2284634939799134
$ws.Range("Q7").Value = 25811.9206573563
$ws.Range("R7").Value = 232307.2859162067
$ws.Range("S7").Value = 0.2201930987898425
$ws.Range("T7").Value = 0.2201930987898425

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("D8").Value = "ECs"
$ws.Range("G8").Value = 2.295346
$ws.Range("H8").Value = 6.886038
$ws.Range("I8").Value = 0.00050911862818281
$ws.Range("J8").Value = 0.0005091186281828102
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.097982
$ws.Range("N8").Value = 15.293946
$ws.Range("O8").Value = 0.1960698660397332
$ws.Range("P8").Value = 0.1960698660397332
$ws.Range("Q8").Value = 11.701632591772
$ws.Range("R8").Value = 105.314693325948
$ws.Range("S8").Value = 0.00009982282122613628
$ws.Range("T8").Value = 0.00009982282122613631

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("D9").Value = "FAPs"
$ws.Range("G9").Value = 2.295346
$ws.Range("H9").Value = 6.886038
$ws.Range("I9").Value = 0.00050911862818281
$ws.Range("J9").Value = 0.0005091186281828102
$ws.Range("M9").Value = 14.96261833333333
$ws.Range("N9").Value = 44.887855
$ws.Range("O9").Value = 0.5754666399803534
$ws.Range("P9").Value = 0.5754666399803534
$ws.Range("Q9").Value = 34.34438614094333
$ws.Range("R9").Value = 309.09947526849
$ws.Range("S9").Value = 0.0002929807863117685
$ws.Range("T9").Value = 0.0002929807863117686

# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.295346
$ws.Range("H10").Value = 6.886038
$ws.Range("I10").Value = 0.00050911862818281
$ws.Range("J10").Value = 0.0005091186281828102
$ws.Range("M10").Value = 5.940243666666667
$ws.Range("N10").Value = 17.820731
$ws.Range("O10").Value = 0.2284634939799134
$ws.Range("P10").Value = 0.2284634939799134
$ws.Range("Q10").Value = 13.63491453930867
$ws.Range("R10").Value = 122.714230853778
$ws.Range("S10").Value = 0.0001163150206449052
$ws.Range("T10").Value = 0.0001163150206449052

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 58.99338900000001
$ws.Range("H11").Value = 176.980167
$ws.Range("I11").Value = 0.01308501344875016
$ws.Range("J11").Value = 0.01308501344875016
$ws.Range("M11").Value = 5.097982
$ws.Range("N11").Value = 15.293946
$ws.Range("O11").Value = 0.1960698660397332
$ws.Range("P11").Value = 0.1960698660397332
$ws.Range("Q11").Value = 300.7472352409981
$ws.Range("R11").Value = 2706.725117168982
$ws.Range("S11").Value = 0.002565576834024551
$ws.Range("T11").Value = 0.002565576834024551

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 58.99338900000001
$ws.Range("H12").Value = 176.980167
$ws.Range("I12").Value = 0.01308501344875016
$ws.Range("J12").Value = 0.01308501344875016
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 14.96261833333333
$ws.Range("N12").Value = 44.887855
$ws.Range("O12").Value = 0.5754666399803534
$ws.Range("P12").Value = 0.5754666399803534
$ws.Range("Q12").Value = 882.6955637968651
$ws.Range("R12").Value = 7944.260074171786
$ws.Range("S12").Value = 0.00752998872344999
$ws.Range("T12").Value = 0.00752998872344999

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 58.99338900000001
$ws.Range("H13").Value = 176.980167
$ws.Range("I13").Value = 0.01308501344875016
$ws.Range("J13").Value = 0.01308501344875016
$ws.Range("M13").Value = 5.940243666666667
$ws.Range("N13").Value = 17.820731
$ws.Range("O13").Value = 0.2284634939799134
$ws.Range("P13").Value = 0.2284634939799134
$ws.Range("Q13").Value = 350.4351053824531
$ws.Range("R13").Value = 3153.915948442078
$ws.Range("S13").Value = 0.002989447891275618
$ws.Range("T13").Value = 0.002989447891275618

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("G14").Value = 2.542663333333333
$ws.Range("H14").Value = 7.62799
$ws.Range("I14").Value = 0.0005639747855867473
$ws.Range("J14").Value = 0.0005639747855867473
$ws.Range("M14").Value = 5.097982
$ws.Range("N14").Value = 15.293946
$ws.Range("O14").Value = 0.1960698660397332
$ws.Range("P14").Value = 0.1960698660397332
$ws.Range("Q14").Value = 12.96245190539333
$ws.Range("R14").Value = 116.66206714854
$ws.Range("S14").Value = 0.0001105784606597808
$ws.Range("T14").Value = 0.0001105784606597808

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("G15").Value = 2.542663333333333
$ws.Range("H15").Value = 7.62799
$ws.Range("I15").Value = 0.0005639747855867473
$ws.Range("J15").Value = 0.0005639747855867473
$ws.Range("O15").Value = 0.5754666399803534
$ws.Range("P15").Value = 0.5754666399803534
$ws.Range("Q15").Value = 38.04490100682778
$ws.Range("R15").Value = 342.40410906145
$ws.Range("S15").Value = 0.0003245486748952457
$ws.Range("T15").Value = 0.0003245486748952457

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("G16").Value = 2.542663333333333
$ws.Range("H16").Value = 7.62799
$ws.Range("I16").Value = 0.0005639747855867473
$ws.Range("J16").Value = 0.0005639747855867473
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 5.940243666666667
$ws.Range("N16").Value = 17.820731
$ws.Range("O16").Value = 0.2284634939799134
$ws.Range("P16").Value = 0.2284634939799134
$ws.Range("Q16").Value = 15.10403976229889
$ws.Range("R16").Value = 135.93635786069
$ws.Range("S16").Value = 0.0001288476500317208
$ws.Range("T16").Value = 0.0001288476500317208

# Remove trailing rows 17-21 (data no longer present in updated export)
$ws.Rows("17:21").Delete()
